$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain-text values that look numeric (e.g. "322.52",
# "29.183.31"). Excel auto-converts a numeric-looking string assigned via
# .Value into a real number, which would corrupt these text labels. Force the
# cells to Text format before writing, then restore the default "Normal" style
# afterwards so no stray formatting is left behind.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.183.31'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.887.57'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '322.52'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4698'
$ws.Range('E7').Value = '  +2.25%  '
$ws.Range('D8').Value = '0.4024'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '47.29'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = '0.08010'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = '0.9929'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = '22.65'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').Value = '1.892.23'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '5.913'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').Value = '7.007'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').Value = '89.21'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '0.06639'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '0.00001024'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = '17.40'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = '29.179.87'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = '5.483'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '11.65'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').Value = '2.175'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').Value = '2.074.15'
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D27').Value = '154.71'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').Value = '19.59'
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('D29').Value = '6.008'
$ws.Range('E29').Value = '  +6.76%  '
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('D31').Value = '117.17'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').Value = '1.022'
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('D33').Value = '0.09414'
$ws.Range('D34').Value = '3.539'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '1.376'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').Value = '5.346'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '0.06044'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D39').Value = '1.169'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = '7.972'
$ws.Range('E40').Value = '  -5.38%  '
$ws.Range('D41').Value = '0.5806'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').Value = '0.1826'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '10.01'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '2.433'
$ws.Range('E44').Value = '  +5.13%  '
$ws.Range('D45').Value = '1.272'
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('D46').Value = '0.07701'
$ws.Range('E46').Value = '  +2.53%  '
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '0.5463'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('D49').Value = '1.899'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('D50').Value = '113.18'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '0.2941'
$ws.Range('E51').Value = '  -0.01%  '

$ws.Range('D2:D51').Style = 'Normal'
